$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (58) down to the
# three new rows (59-61) so the new cells pick up the same styles
# (bold/bordered index column, datetime-formatted date column, etc.)
$ws.Range("A58:V58").Copy()
$ws.Range("A59:V61").PasteSpecial(-4122)

# ---- Row 59 (Indice 58) ----
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "argentina"
$ws.Range("C59").Value = "copa-de-la-liga-profesional"
$ws.Range("D59").Value = "2023"
$ws.Range("E59").Value = 45188.98958333334
$ws.Range("F59").Value = "Central Cordoba"
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = "Boca Juniors"
$ws.Range("I59").Value = 3
$ws.Range("J59").Value = 3.44
$ws.Range("K59").Value = "15/09/2023 23:13"
$ws.Range("L59").Value = 3.57
$ws.Range("M59").Value = "19/09/2023 23:43"
$ws.Range("N59").Value = 2.98
$ws.Range("O59").Value = "15/09/2023 23:13"
$ws.Range("P59").Value = 3.12
$ws.Range("Q59").Value = "19/09/2023 23:43"
$ws.Range("R59").Value = 2.4
$ws.Range("S59").Value = "15/09/2023 23:13"
$ws.Range("T59").Value = 2.31
$ws.Range("U59").Value = "19/09/2023 23:43"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/central-cordoba-santiago-del-estero-boca-juniors/CYd8EYPG/"

# ---- Row 60 (Indice 59) ----
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "argentina"
$ws.Range("C60").Value = "copa-de-la-liga-profesional"
$ws.Range("D60").Value = "2023"
$ws.Range("E60").Value = 45189.08333333334
$ws.Range("F60").Value = "Huracan"
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = "Gimnasia L.P."
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1.8
$ws.Range("K60").Value = "15/09/2023 01:13"
$ws.Range("L60").Value = 1.9
$ws.Range("M60").Value = "20/09/2023 01:55"
$ws.Range("N60").Value = 3.3
$ws.Range("O60").Value = "15/09/2023 01:13"
$ws.Range("P60").Value = 3.1
$ws.Range("Q60").Value = "20/09/2023 01:55"
$ws.Range("R60").Value = 4.82
$ws.Range("S60").Value = "15/09/2023 01:13"
$ws.Range("T60").Value = 5.38
$ws.Range("U60").Value = "20/09/2023 01:55"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/huracan-gimnasia-l-p/IZ3aGfe4/"

# ---- Row 61 (Indice 60) ----
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "argentina"
$ws.Range("C61").Value = "copa-de-la-liga-profesional"
$ws.Range("D61").Value = "2023"
$ws.Range("E61").Value = 45189.08333333334
$ws.Range("F61").Value = "Rosario Central"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = "Independiente"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 2.1
$ws.Range("K61").Value = "16/09/2023 01:13"
$ws.Range("L61").Value = 2.78
$ws.Range("M61").Value = "20/09/2023 01:55"
$ws.Range("N61").Value = 3.12
$ws.Range("O61").Value = "16/09/2023 01:13"
$ws.Range("P61").Value = 3.02
$ws.Range("Q61").Value = "20/09/2023 01:55"
$ws.Range("R61").Value = 4.04
$ws.Range("S61").Value = "16/09/2023 01:13"
$ws.Range("T61").Value = 2.93
$ws.Range("U61").Value = "20/09/2023 01:55"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/rosario-central-independiente/vq5eHztb/"

# The "temporada" column (D) stores the year as TEXT in the source data
# ("2023" as a string), not a number. A direct .Value assignment of the
# string "2023" gets auto-converted to a numeric value by Excel. D2
# already holds that same text "2023" with the default (unstyled)
# formatting, so reuse it as the source for a values-only paste into the
# new D cells - this keeps them as text without introducing any new
# number formats/styles.
$ws.Range("D2").Copy()
$ws.Range("D59").PasteSpecial(-4163)
$ws.Range("D2").Copy()
$ws.Range("D60").PasteSpecial(-4163)
$ws.Range("D2").Copy()
$ws.Range("D61").PasteSpecial(-4163)
